$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 50) with the latest bitcoin buy entry (run on 2025-09-28).
# Column A holds the date as plain text (matching the existing rows that were
# appended as text dates), so force a text format before assigning the value
# to stop Excel from auto-converting the "mm/dd/yyyy"-looking string into a
# date serial number, then clear the formatting so the cell ends up with no
# explicit style, just like the other text-date rows in this sheet.
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "09/28/2025"
$ws.Range("A50").ClearFormats()

$ws.Range("B50").Value = 0.0004548000000000017
$ws.Range("C50").Value = 109938.4344766926
$ws.Range("D50").Value = 50
